$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last_update" timestamp in B2 (footer live-updated info)
$ws.Range("B2").Value = 1728846186

# Set explicit width on column B (bestFit-equivalent width captured by Excel)
$ws.Columns("B").ColumnWidth = 10.33

# Move/restore the active selection to D8, matching the saved view state
$ws.Range("D8").Select()
